$d = $word.ActiveDocument

# 1) Merge "Zalozona w 20" + "18" + " roku Los " runs into a single run.
$d.Content.Find.Execute("Założona w 2018 roku Los ", $true, $false, $false, $false, $false, $true, 1, $false, "Założona w 2018 roku Los ", 2) | Out-Null

# 2) Merge "Dzial IT ... - zapewniaja ..." runs into a single run.
$d.Content.Find.Execute("Dział IT – odpowiada za utrzymanie i rozwój systemów informatycznych, które - zapewniają sprawne funkcjonowanie sieci logistycznej.", $true, $false, $false, $false, $false, $true, 1, $false, "Dział IT – odpowiada za utrzymanie i rozwój systemów informatycznych, które - zapewniają sprawne funkcjonowanie sieci logistycznej.", 2) | Out-Null

# 3) Shorten the figure caption text above the diagram.
$d.Content.Find.Execute("Kontekstowy diagram przypadków użycia (DPU) projektowanego systemu", $true, $false, $false, $false, $false, $true, 1, $false, "Kontekstowy diagram projektowanego systemu", 2) | Out-Null

# 4) Rename "Biznesowy kontekst" to "Diagram kontekstowy" in the bold diagram title.
$d.Content.Find.Execute("Diagram 1. Biznesowy kontekst Systemu Zarządzania Przesyłkami", $true, $false, $false, $false, $false, $true, 1, $false, "Diagram 1. Diagram kontekstowy Systemu Zarządzania Przesyłkami", 2) | Out-Null
